$d = $word.ActiveDocument

$replacements = @(
    @{Old = "2024-10-02 Wednesday"; New = "2024-10-03 Thursday"},
    @{Old = "55×66="; New = "15×61="},
    @{Old = "95×86="; New = "53×91="},
    @{Old = "68×93="; New = "78×44="},
    @{Old = "54×51="; New = "69×58="},
    @{Old = "30×32="; New = "11×67="},
    @{Old = "21×49="; New = "32×21="},
    @{Old = "81×66="; New = "25×68="},
    @{Old = "72×94="; New = "14×57="},
    @{Old = "27×38="; New = "11×42="},
    @{Old = "94×55="; New = "37×53="},
    @{Old = "82×11="; New = "51×92="},
    @{Old = "31×36="; New = "37×84="},
    @{Old = "58×74="; New = "96×92="},
    @{Old = "56×74="; New = "15×68="},
    @{Old = "41×33="; New = "33×50="},
    @{Old = "29×83="; New = "88×52="},
    @{Old = "64×94="; New = "85×96="},
    @{Old = "65×40="; New = "64×13="},
    @{Old = "44×56="; New = "33×17="},
    @{Old = "42×68="; New = "46×49="},
    @{Old = "66×98="; New = "80×65="},
    @{Old = "58×32="; New = "59×25="},
    @{Old = "77×73="; New = "68×60="},
    @{Old = "76×93="; New = "54×47="},
    @{Old = "69×27="; New = "37×74="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

$d.Save()
